$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 39; this shifts the existing rows 39-55 down
# to 40-56, preserving all of their data and formatting.
$ws.Rows("39:39").Insert()

# Populate the newly inserted row 39 with the new weekly price record.
$ws.Range("A39").Value = 10
$ws.Range("B39").Value = "Vega Modelo de Temuco"
$ws.Range("C39").Value = "La Araucanía"
$ws.Range("D39").Value = 44777
$ws.Range("E39").Value = 9
$ws.Range("F39").Value = "Fruta"
$ws.Range("G39").Value = 100108
$ws.Range("H39").Value = "Tropicales y subtropicales"
$ws.Range("I39").Value = 100108003
$ws.Range("J39").Value = "Maracuyá"
$ws.Range("K39").Value = "Sin especificar"
$ws.Range("L39").Value = "Primera"
$ws.Range("M39").Value = 50
$ws.Range("N39").Value = 35000
$ws.Range("O39").Value = 35000
$ws.Range("P39").Value = 35000
$ws.Range("Q39").Value = "$/caja 18 kilos"
$ws.Range("R39").Value = "Región de Arica y Parinacota"
$ws.Range("S39").Value = 1944
$ws.Range("T39").Value = 18
